$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the leading gap (rows 3-7, currently empty); everything
#    below shifts up by 5, sliding the stats block from rows 8-17 into
#    rows 3-12, and the trailing blank row 18 up to row 13. Because
#    it's a row delete (not a plain cell move), internal relative
#    references (e.g. "B8" inside the skewness formula) are translated
#    along with it automatically.
$ws.Rows("3:7").Delete()

# The engine re-derives a couple of the single-value stat formulas with
# an extra coercion wrapper on row delete; restore their exact text.
$ws.Range("B4").Formula = "=_xlfn.VAR.S(A2:AP2)"
$ws.Range("B5").Formula = "=_xlfn.STDEV.S(A2:AP2)"
$ws.Range("B6").Formula = "=_xlfn.MODE.SNGL(A2:AP2)"

# 2) The trailing blank row (now at row 13) belongs back at row 18, with
#    nothing in between. Insert 5 fresh rows above it -- this carries
#    its existing formatted-but-empty cells down to row 18 intact --
#    then fully clear the 5 newly inserted placeholder rows so rows
#    13-17 are genuinely empty (no leftover formatting).
$ws.Rows("13:17").Insert()
$ws.Range("A13:B17").Clear()

# Update cursor/selection to match the post-edit state.
$ws.Range("G12").Select()
